$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora): all rows 2-51 change from 16 to 17
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "17"

# Column D (Price) updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.73"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.81"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08064"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.914"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.837"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9299"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1322"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1904"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09216"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03476"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09897"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001419"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006628"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.609"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.014"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.164"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2535"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04418"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004710"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003131"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01992"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05216"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007613"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01015"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1366"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01073"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006318"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001601"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"

# Column E (Volume 1h) updates
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.12%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.47%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.15%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.32%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.86%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.01%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.74%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.23%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.16%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.25%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.73%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.38%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.31%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.96%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "11.86%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.30%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.57%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.31%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.75%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.73%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.72%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.21%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.77%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.36%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.53%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "9.51%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.55%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.44%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.16%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.40%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.07%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.24%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.47%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.49%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
